# Insert one new data row at row 399 (pushing the existing rows 399-435
# down to 400-436) in the "Hortaliza, Terminal La Palmera de La Serena -
# Espinaca" weekly price sheet, and populate it with a new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 399:435 down one row, carrying the row-399 formatting
# (this is what gives new D399 the date style used throughout column D).
$ws.Rows("399:399").Insert()

# Fill in the newly inserted row with the new price observation.
$ws.Range("A399").Value = 8
$ws.Range("B399").Value = "Terminal La Palmera de La Serena"
$ws.Range("C399").Value = "Coquimbo"
$ws.Range("D399").Value = 45106
$ws.Range("E399").Value = 4
$ws.Range("F399").Value = 100112012
$ws.Range("G399").Value = "Espinaca"
$ws.Range("H399").Value = "Sin especificar"
$ws.Range("I399").Value = "Primera"
$ws.Range("J399").Value = 1000
$ws.Range("K399").Value = 400
$ws.Range("L399").Value = 500
$ws.Range("M399").Value = 450
$ws.Range("N399").Value = "$/atado 300 a 500 gramos"
$ws.Range("O399").Value = "Provincia del Elquí"
$ws.Range("P399").Value = 900
$ws.Range("Q399").Value = 0.5
$ws.Range("R399").Value = "Hortaliza"
